$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 600.7273
$ws.Range("I11").Value = 600.7273
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 600.7273
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -460.7273
$ws.Range("H38").Value = 585.3333
$ws.Range("I38").Value = 283.5
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 850.5
$ws.Range("L38").Value = 9000
$ws.Range("M38").Value = -478.5
$ws.Range("N38").Value = -9744
$ws.Range("H39").Value = 682.55554
$ws.Range("I39").Value = 144.125
$ws.Range("J39").Value = 4990
$ws.Range("K39").Value = 432.375
$ws.Range("L39").Value = 14970
$ws.Range("M39").Value = -136.375
$ws.Range("N39").Value = -15562
$ws.Range("H74").Value = 3666.6667
$ws.Range("I74").Value = 3500
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 3500
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -2564
$ws.Range("N74").Value = -5872
$ws.Range("H77").Value = 3666.6667
$ws.Range("I77").Value = 3500
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 17500
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -12820
$ws.Range("N77").Value = -29360
$ws.Range("H100").Value = 5299.6
$ws.Range("I100").Value = 5299.6
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 5299.6
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = None
$ws.Range("N100").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H37").Value = 28500
$ws.Range("I37").Value = 23750
$ws.Range("J37").Value = 38000
$ws.Range("K37").Value = 23750
$ws.Range("L37").Value = 38000
$ws.Range("M37").Value = -23477
$ws.Range("N37").Value = -38546
$ws.Range("H61").Value = 6832
$ws.Range("I61").Value = 6832
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 6832
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -6620
$ws.Range("H74").Value = 1468.2222
$ws.Range("I74").Value = 901.875
$ws.Range("J74").Value = 5999
$ws.Range("K74").Value = 901.875
$ws.Range("L74").Value = 5999
$ws.Range("M74").Value = -27.875
$ws.Range("N74").Value = -7747
$ws.Range("H77").Value = 1468.2222
$ws.Range("I77").Value = 901.875
$ws.Range("J77").Value = 5999
$ws.Range("K77").Value = 4509.375
$ws.Range("L77").Value = 29995
$ws.Range("M77").Value = -141.375
$ws.Range("N77").Value = -38731
$ws.Range("H92").Value = 43883.332
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 43883.332
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 43883.332
$ws.Range("N92").Value = -48875.332
$ws.Range("H132").Value = 1546.125
$ws.Range("I132").Value = 1378.9231
$ws.Range("J132").Value = 2270.6667
$ws.Range("K132").Value = 4136.7693
$ws.Range("L132").Value = 6812.000100000001
$ws.Range("M132").Value = -1606.7693
$ws.Range("N132").Value = -11872.0001
$ws.Range("H136").Value = 6832
$ws.Range("I136").Value = 6832
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 20496
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -17946

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3698.75
$ws.Range("I86").Value = 3598.3333
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 3598.3333
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -2475.3333
$ws.Range("N86").Value = -6246
$ws.Range("H89").Value = 3698.75
$ws.Range("I89").Value = 3598.3333
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 17991.6665
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -12375.6665
$ws.Range("N89").Value = -31232
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = None
$ws.Range("N92").ClearContents()
$ws.Range("H105").Value = 4518.5
$ws.Range("I105").Value = 4465
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 4465
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -2718
$ws.Range("N105").Value = -8494
$ws.Range("H131").Value = 70000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 70000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 70000
$ws.Range("N131").Value = -80080

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 8070.7144
$ws.Range("I41").Value = 3582.5
$ws.Range("J41").Value = 35000
$ws.Range("K41").Value = 3582.5
$ws.Range("L41").Value = 35000
$ws.Range("M41").Value = -3154.5
$ws.Range("N41").Value = -35856
$ws.Range("H50").Value = 26000.334
$ws.Range("I50").Value = 10000
$ws.Range("J50").Value = 34000.5
$ws.Range("K50").Value = 10000
$ws.Range("L50").Value = 34000.5
$ws.Range("M50").Value = -9375
$ws.Range("N50").Value = -35250.5
$ws.Range("H51").Value = 32500
$ws.Range("I51").Value = 25000
$ws.Range("J51").Value = 40000
$ws.Range("K51").Value = 25000
$ws.Range("L51").Value = 40000
$ws.Range("M51").Value = -24264
$ws.Range("N51").Value = -41472
$ws.Range("H58").Value = 4125.5
$ws.Range("I58").Value = 4125.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 4125.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -3922.5
$ws.Range("H59").Value = 44665.668
$ws.Range("I59").Value = 41998
$ws.Range("J59").Value = 45999.5
$ws.Range("K59").Value = 41998
$ws.Range("L59").Value = 45999.5
$ws.Range("M59").Value = -40853
$ws.Range("N59").Value = -48289.5
$ws.Range("H60").Value = 31500
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 31500
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = None
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -32522
$ws.Range("H61").Value = 32500
$ws.Range("I61").Value = 25000
$ws.Range("J61").Value = 40000
$ws.Range("K61").Value = 25000
$ws.Range("L61").Value = 40000
$ws.Range("M61").Value = -24652
$ws.Range("N61").Value = -40696
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = None
$ws.Range("N131").ClearContents()
$ws.Range("H136").Value = 4125.5
$ws.Range("I136").Value = 4125.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 12376.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -9826.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 159.875
$ws.Range("I6").Value = 168.57143
$ws.Range("J6").Value = 99
$ws.Range("K6").Value = 505.71429
$ws.Range("L6").Value = 297
$ws.Range("M6").Value = -392.71429
$ws.Range("N6").Value = -523
$ws.Range("H140").Value = 2379.25
$ws.Range("I140").Value = 811
$ws.Range("J140").Value = 4993
$ws.Range("K140").Value = 2433
$ws.Range("L140").Value = 14979
$ws.Range("M140").Value = 2747
$ws.Range("N140").Value = -25339

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3247.4
$ws.Range("I43").Value = 2119.5
$ws.Range("J43").Value = 3999.3333
$ws.Range("K43").Value = 2119.5
$ws.Range("L43").Value = 3999.3333
$ws.Range("M43").Value = -1968.5
$ws.Range("N43").Value = -4301.3333
$ws.Range("H46").Value = 9268.9
$ws.Range("I46").Value = 3480
$ws.Range("J46").Value = 15057.8
$ws.Range("K46").Value = 3480
$ws.Range("L46").Value = 15057.8
$ws.Range("M46").Value = -3324
$ws.Range("N46").Value = -15369.8
$ws.Range("H57").Value = 19998.75
$ws.Range("I57").Value = 20000
$ws.Range("J57").Value = 19997.5
$ws.Range("K57").Value = 20000
$ws.Range("L57").Value = 19997.5
$ws.Range("M57").Value = -19180
$ws.Range("N57").Value = -21637.5
$ws.Range("H70").Value = 12799.667
$ws.Range("I70").Value = 12799.667
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 12799.667
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -12529.667
$ws.Range("H73").Value = 12799.667
$ws.Range("I73").Value = 12799.667
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 12799.667
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -11863.667
$ws.Range("H80").Value = 1699.8
$ws.Range("I80").Value = 1249.5
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 1249.5
$ws.Range("L80").Value = 2000
$ws.Range("M80").Value = -251.5
$ws.Range("N80").Value = -3996
$ws.Range("H83").Value = 1699.8
$ws.Range("I83").Value = 1249.5
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 6247.5
$ws.Range("L83").Value = 10000
$ws.Range("M83").Value = -1255.5
$ws.Range("N83").Value = -19984
$ws.Range("H92").Value = 15299.333
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 15299.333
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 15299.333
$ws.Range("N92").Value = -19043.333

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2257
$ws.Range("I22").Value = 2349
$ws.Range("J22").Value = 1797
$ws.Range("K22").Value = 2349
$ws.Range("L22").Value = 1797
$ws.Range("M22").Value = -2054
$ws.Range("N22").Value = -2387
$ws.Range("H27").Value = 2257
$ws.Range("I27").Value = 2349
$ws.Range("J27").Value = 1797
$ws.Range("K27").Value = 2349
$ws.Range("L27").Value = 1797
$ws.Range("M27").Value = -2242
$ws.Range("N27").Value = -2011
$ws.Range("H136").Value = 3000
$ws.Range("I136").Value = 3000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -6450

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 19997
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 19997
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 19997
$ws.Range("N31").Value = -20693
$ws.Range("H81").Value = 3000.4
$ws.Range("I81").Value = 1667.6666
$ws.Range("J81").Value = 4999.5
$ws.Range("K81").Value = 3335.3332
$ws.Range("L81").Value = 9999
$ws.Range("M81").Value = -2274.3332
$ws.Range("N81").Value = -12121
$ws.Range("H84").Value = 3000.4
$ws.Range("I84").Value = 1667.6666
$ws.Range("J84").Value = 4999.5
$ws.Range("K84").Value = 16676.666
$ws.Range("L84").Value = 49995
$ws.Range("M84").Value = -11372.666
$ws.Range("N84").Value = -60603
$ws.Range("H107").Value = 563.4
$ws.Range("I107").Value = 516.75
$ws.Range("J107").Value = 750
$ws.Range("K107").Value = 1550.25
$ws.Range("L107").Value = 2250
$ws.Range("M107").Value = 369.75
$ws.Range("N107").Value = -6090
$ws.Range("H133").Value = 49997
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 49997
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 49997
$ws.Range("N133").Value = -60117
